$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: January / Squirrel, Male Count 2, Unknown Sex Count 3
$ws.Range("F2").Value = "Squirrel"
$ws.Range("H2").Value = 2
$ws.Range("J2").Value = 3

# Row 3: January / Weasel, Male Count 2, Unknown Sex Count 0
$ws.Range("F3").Value = "Weasel"
$ws.Range("H3").Value = 2
$ws.Range("J3").Value = 0

# Row 4: February / (no species)
$ws.Range("E4").Value = "February"
$ws.Range("F4").Value = ""

# Remove the old row 5 (was the February placeholder row, now redundant)
$ws.Rows("5").Delete()
